$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: convert B2:H2 from text-like values ("5000.00", "20.0", "0.00",
# "100.00") to real numeric values (5000, 20, 0, 100)
$ws.Range("B2").Value = 5000
$ws.Range("C2").Value = 5000
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 5000
$ws.Range("F2").Value = 20
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 100

# Row 3: add a new row of data below, keeping the values as text (so the
# trailing zeros / decimal formatting from the source are preserved
# exactly, instead of being collapsed by automatic number conversion).
$ws.Range("A3:H3").NumberFormat = "@"
$ws.Range("A3").Value = "05/08/2023"
$ws.Range("B3").Value = "4000.00"
$ws.Range("C3").Value = "9000.00"
$ws.Range("D3").Value = "4000.00"
$ws.Range("E3").Value = "9000.00"
$ws.Range("F3").Value = "70.0"
$ws.Range("G3").Value = "0.00"
$ws.Range("H3").Value = "100.00"
# Restore the default (unstyled) cell style for the new row so only the
# data type/content changed, not its formatting.
$ws.Range("A3:H3").Style = "Normal"
